$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 14
$ws.Range("B2").Value = 20
$ws.Range("B3").Value = 14
$ws.Range("B4").Value = 12
$ws.Range("B5").Value = 15

# Restore the default selection (A1) - the saved workbook previously had
# an explicit selection parked on A13 (left over from editing); move it
# back to the top-left cell like a freshly opened/saved workbook.
$ws.Range("A1").Select()
